# Add Q3-2022 data: insert a new "2022-Q3" sheet right after "总计",
# and add the corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new worksheet directly after "总计" (before "2022-Q2")
# ---------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($existingQ2)
$newSheet.Name = "2022-Q3"

# Apply the same "header / index column" styling used throughout the
# workbook: bold font, thin box border, centered horizontal, top vertical.
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------
# 2. Populate header row
# ---------------------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2   # headers start at column B
    $dstCell = $newSheet.Cells.Item(1, $col)
    $dstCell.Value = $headers[$i]
    Set-HeaderStyle $dstCell
}

# ---------------------------------------------------------------
# 3. Populate data rows
# ---------------------------------------------------------------
$rows = @(
    @{A=0;  B="012428"; C="华夏核心制造混合A";                   D="34.97"; E="85.52"; F="4.70"; G="1.6436"; H=7},
    @{A=1;  B="050009"; C="博时新兴成长混合";                     D="24.16"; E="92.80"; F="4.78"; G="1.1548"; H=6},
    @{A=2;  B="013389"; C="华夏成长先锋一年持有混合A";             D="13.25"; E="89.79"; F="4.21"; G="0.5578"; H=9},
    @{A=3;  B="012429"; C="华夏核心制造混合C";                    D="8.86";  E="85.52"; F="4.70"; G="0.4164"; H=7},
    @{A=4;  B="169103"; C="东方红睿轩三年定期开放灵活配置混合";     D="11.31"; E="70.03"; F="2.20"; G="0.2488"; H=10},
    @{A=5;  B="013390"; C="华夏成长先锋一年持有混合C";             D="4.42";  E="89.79"; F="4.21"; G="0.1861"; H=9},
    @{A=6;  B="159743"; C="博时中证湖北新旧动能转换ETF";           D="3.31";  E="98.93"; F="1.78"; G="0.0589"; H=7},
    @{A=7;  B="014737"; C="创金合信专精特新股票C";                D="0.58";  E="81.21"; F="6.12"; G="0.0355"; H=4},
    @{A=8;  B="002863"; C="金信深圳成长灵活配置混合";              D="0.61";  E="93.95"; F="5.19"; G="0.0317"; H=5},
    @{A=9;  B="014736"; C="创金合信专精特新股票A";                D="0.35";  E="81.21"; F="6.12"; G="0.0214"; H=4},
    @{A=10; B="011377"; C="创金合信积极成长股票A";                D="0.18";  E="93.65"; F="8.45"; G="0.0152"; H=5},
    @{A=11; B="011378"; C="创金合信积极成长股票C";                D="0.13";  E="93.65"; F="8.45"; G="0.0110"; H=5}
)

$r = 2
foreach ($row in $rows) {
    $aCell = $newSheet.Cells.Item($r, 1)
    $aCell.Value = $row.A
    Set-HeaderStyle $aCell

    $newSheet.Cells.Item($r, 2).Value = $row.B
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row.B

    $newSheet.Cells.Item($r, 3).Value = $row.C

    foreach ($col in 4..7) {
        $cell = $newSheet.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
    }
    $newSheet.Cells.Item($r, 4).Value = $row.D
    $newSheet.Cells.Item($r, 5).Value = $row.E
    $newSheet.Cells.Item($r, 6).Value = $row.F
    $newSheet.Cells.Item($r, 7).Value = $row.G

    $newSheet.Cells.Item($r, 8).Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert the new 2022-Q3 row
#    and shift the existing rows down by one.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryRows = @(
    @{A=0; B="2022-Q3"; C=12; D=4.38},
    @{A=1; B="2022-Q2"; C=6;  D=4.33},
    @{A=2; B="2022-Q1"; C=7;  D=4.98},
    @{A=3; B="2021-Q4"; C=32; D=15.66},
    @{A=4; B="2021-Q3"; C=7;  D=5.41},
    @{A=5; B="2021-Q2"; C=12; D=4.34},
    @{A=6; B="2021-Q1"; C=16; D=4.38},
    @{A=7; B="2020-Q4"; C=10; D=3.32}
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row.A
    $summary.Cells.Item($r, 2).Value = $row.B
    $summary.Cells.Item($r, 3).Value = $row.C
    $summary.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}

# Restore "总计" as the active sheet (Worksheets.Add left the newly
# inserted sheet active).
$summary.Activate()
